$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header relabel: "Name" -> "File_info", "File 1"/"File 2" -> "File"/"File"
$ws.Range("A1").Value = "File_info"
$ws.Range("B1").Value = "File"
$ws.Range("C1").Value = "File"

# Row 2 label relabel: "File" -> "Filename" (B2/C2 unchanged)
$ws.Range("A2").Value = "Filename"

# Update the active selection to C1 (was A2)
$ws.Range("C1").Select()
